$wb = $excel.ActiveWorkbook

# "Metadata" sheet holds property/value pairs.
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 = Experimental -> set value to the literal text "false"
# (a plain .Value assignment of "false" is auto-coerced to a real Boolean
# cell by Excel; route it through a text formula + paste-as-values so the
# result lands as literal text, matching the source data.)
$ws.Cells.Item(7, 2).Formula = "=""fal""&""se"""
$ws.Cells.Item(7, 2).Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4163)

# Row 8 = Date -> update timestamp
$ws.Range("B8").Value = "2025-11-04T10:04:56+00:00"
